# Auto-generated edit script: refresh crypto price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.531.52'
$ws.Range('E2').Value = '  +1.07%  '
$ws.Range('D3').Value = '1.850.55'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.42'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4718'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.75%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2738'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06309'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '17.66'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.75%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07448'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.49%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.779.18'
$ws.Range('E12').Value = '  -3.85%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.998'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.14%  '
$ws.Range('E14').Value = '  -0.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6242'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.50%  '
$ws.Range('D16').Value = '30.486.18'
$ws.Range('E16').Value = '  +1.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '242.97'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +6.94%  '
$ws.Range('E18').Value = '  -0.01%  '
$ws.Range('E19').Value = '  +0.51%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007322'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9992'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.08%  '
$ws.Range('E22').Value = '  -0.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.916'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.153'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.86%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '162.61'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '17.92'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.93%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.876'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.51%  '
$ws.Range('E28').Value = '  -1.48%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.360'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.78%  '
$ws.Range('E30').Value = '  -2.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.816'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.04844'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.90%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.132'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.05%  '
$ws.Range('E34').Value = '  -1.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.708'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.38%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.01895'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.26%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.676'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.52%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.8731'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.44%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.987'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.82%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '106.57'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.50%  '
$ws.Range('E41').Value = '  +0.23%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.522'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.31%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4053'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.162'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.62%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '62.40'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.76%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1210'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '33.42'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.99%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.504'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05542'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.352'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3655'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.31%  '
